$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M18").Value = -164.22223
$ws.Range("H18").Value = 448.22223
$ws.Range("K18").Value = 448.22223
$ws.Range("I18").Value = 448.22223
$ws.Range("J32").Value = 6480.5
$ws.Range("H32").Value = 5297.875
$ws.Range("N32").Value = -7132.5
$ws.Range("L32").Value = 6480.5
$ws.Range("K101").Value = 1884
$ws.Range("M101").Value = -262
$ws.Range("H101").Value = 831.6667
$ws.Range("I101").Value = 628
$ws.Range("K137").Value = 145079.145
$ws.Range("L137").Value = 18866.334
$ws.Range("H137").Value = 24694.812
$ws.Range("I137").Value = 48359.715
$ws.Range("N137").Value = -23966.334
$ws.Range("M137").Value = -142529.145
$ws.Range("J137").Value = 6288.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 3269.0518
$ws.Range("I32").Value = 3269.0518
$ws.Range("H32").Value = 3164.2856
$ws.Range("M32").Value = -2982.0518
$ws.Range("H74").Value = 37444.152
$ws.Range("I74").Value = 38488.844
$ws.Range("N74").Value = -5762
$ws.Range("J74").Value = 4014
$ws.Range("M74").Value = -37614.844
$ws.Range("L74").Value = 4014
$ws.Range("K74").Value = 38488.844
$ws.Range("L77").Value = 20070
$ws.Range("M77").Value = -188076.22
$ws.Range("K77").Value = 192444.22
$ws.Range("I77").Value = 38488.844
$ws.Range("J77").Value = 4014
$ws.Range("N77").Value = -28806
$ws.Range("H77").Value = 37444.152
$ws.Range("H132").Value = 55013.63
$ws.Range("M132").Value = -6336.2855
$ws.Range("K132").Value = 8866.2855
$ws.Range("I132").Value = 2955.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1309.5
$ws.Range("L22").Value = 1869
$ws.Range("J22").Value = 1869
$ws.Range("N22").Value = -2215
$ws.Range("H86").Value = 51008.332
$ws.Range("I86").Value = 26087.5
$ws.Range("K86").Value = 26087.5
$ws.Range("M86").Value = -24964.5
$ws.Range("M89").Value = -124821.5
$ws.Range("K89").Value = 130437.5
$ws.Range("H89").Value = 51008.332
$ws.Range("I89").Value = 26087.5
$ws.Range("M94").Value = -342.7368
$ws.Range("H94").Value = 1368.6061
$ws.Range("K94").Value = 793.7368
$ws.Range("I94").Value = 793.7368
$ws.Range("L123").Value = 85000
$ws.Range("N123").Value = -94800
$ws.Range("H123").Value = 85000
$ws.Range("J123").Value = 85000
$ws.Range("H132").Value = 118991.664
$ws.Range("N132").Value = -129111.664
$ws.Range("J132").Value = 118991.664
$ws.Range("L132").Value = 118991.664
$ws.Range("K134").Value = 13673.4
$ws.Range("M134").Value = -11138.4
$ws.Range("L134").Value = 27042
$ws.Range("H134").Value = 5300.5
$ws.Range("J134").Value = 9014
$ws.Range("I134").Value = 4557.8
$ws.Range("N134").Value = -32112

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("H86").Value = 9990
$ws.Range("I86").Value = 9990
$ws.Range("N86").Value = 0
$ws.Range("K86").Value = 9990
$ws.Range("M86").Value = -8867
$ws.Range("M89").Value = -44334
$ws.Range("N89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 49950
$ws.Range("L89").ClearContents()
$ws.Range("H89").Value = 9990
$ws.Range("I89").Value = 9990
$ws.Range("H99").Value = 3599.5
$ws.Range("I99").Value = 3466.6667
$ws.Range("K99").Value = 3466.6667
$ws.Range("M99").Value = -1968.6667
$ws.Range("K122").Value = 5788.7142
$ws.Range("L122").Value = 10722.75
$ws.Range("M122").Value = -3338.7142
$ws.Range("N122").Value = -15622.75
$ws.Range("H122").Value = 2527.6365
$ws.Range("J122").Value = 3574.25
$ws.Range("I122").Value = 1929.5714
$ws.Range("H126").Value = 3599.5
$ws.Range("I126").Value = 3466.6667
$ws.Range("K126").Value = 10400.0001
$ws.Range("M126").Value = -7930.000100000001
$ws.Range("H132").Value = 3390.054
$ws.Range("M132").Value = -7499
$ws.Range("L132").Value = 11334.75
$ws.Range("J132").Value = 3778.25
$ws.Range("K132").Value = 10029
$ws.Range("N132").Value = -16394.75
$ws.Range("I132").Value = 3343
$ws.Range("K134").Value = 3738.375
$ws.Range("M134").Value = -1203.375
$ws.Range("L134").Value = 6554.499899999999
$ws.Range("H134").Value = 1648.4286
$ws.Range("J134").Value = 2184.8333
$ws.Range("I134").Value = 1246.125
$ws.Range("N134").Value = -11624.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M9").Value = -2626
$ws.Range("K9").Value = 2850
$ws.Range("H9").Value = 950
$ws.Range("I9").Value = 950
$ws.Range("K55").Value = 13500300
$ws.Range("I55").Value = 4500100
$ws.Range("M55").Value = -13500123
$ws.Range("H55").Value = 4500100
$ws.Range("J107").Value = 3832.6667
$ws.Range("L107").Value = 11498.0001
$ws.Range("H107").Value = 2154
$ws.Range("N107").Value = -15338.0001
$ws.Range("K107").Value = 1426.00002
$ws.Range("I107").Value = 475.33334
$ws.Range("M107").Value = 493.9999800000001
$ws.Range("K120").Value = 11175
$ws.Range("M120").Value = -6337
$ws.Range("I120").Value = 3725
$ws.Range("H120").Value = 3725
$ws.Range("J131").Value = 1777.6666
$ws.Range("N131").Value = -15412.9998
$ws.Range("H131").Value = 21740510
$ws.Range("L131").Value = 5332.9998
$ws.Range("L133").Value = 14992.5
$ws.Range("N133").Value = -25112.5
$ws.Range("J133").Value = 4997.5
$ws.Range("I133").Value = 1776.6666
$ws.Range("H133").Value = 3065
$ws.Range("M133").Value = -269.9997999999996
$ws.Range("K133").Value = 5329.9998
$ws.Range("H134").Value = 1289.3334
$ws.Range("L134").Value = 10340.0001
$ws.Range("N134").Value = -20480.0001
$ws.Range("J134").Value = 3446.6667
$ws.Range("K136").Value = 4229.625
$ws.Range("M136").Value = 870.375
$ws.Range("I136").Value = 1409.875
$ws.Range("H136").Value = 1409.875
$ws.Range("K137").Value = 7038
$ws.Range("L137").Value = 9099
$ws.Range("H137").Value = 2460.5
$ws.Range("I137").Value = 2346
$ws.Range("N137").Value = -19299
$ws.Range("M137").Value = -1938
$ws.Range("J137").Value = 3033

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L102").Value = 8139.9
$ws.Range("N102").Value = -11383.9
$ws.Range("M102").Value = -756.875
$ws.Range("K102").Value = 2378.875
$ws.Range("I102").Value = 2378.875
$ws.Range("H102").Value = 4594.654
$ws.Range("J102").Value = 8139.9
$ws.Range("H122").Value = 3505.8
$ws.Range("J122").Value = 5250
$ws.Range("L122").Value = 15750
$ws.Range("N122").Value = -20650
$ws.Range("H126").Value = 3445.0833
$ws.Range("I126").Value = 2342.2856
$ws.Range("K126").Value = 7026.8568
$ws.Range("M126").Value = -4556.8568
$ws.Range("H132").Value = 1962.4517
$ws.Range("M132").Value = -1546.25
$ws.Range("K132").Value = 4076.25
$ws.Range("I132").Value = 1358.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 4000
$ws.Range("M40").Value = -3864
$ws.Range("H40").Value = 7997.25
$ws.Range("K40").Value = 4000
$ws.Range("H46").Value = 2124.875
$ws.Range("H74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("K122").Value = 12256.7724
$ws.Range("L122").Value = 15110.25
$ws.Range("M122").Value = -9806.7724
$ws.Range("N122").Value = -20010.25
$ws.Range("H122").Value = 4339.2334
$ws.Range("J122").Value = 5036.75
$ws.Range("I122").Value = 4085.5908
$ws.Range("L123").Value = 74996
$ws.Range("N123").Value = -84796
$ws.Range("H123").Value = 74996
$ws.Range("J123").Value = 74996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 44999.75
$ws.Range("J75").Value = 44999.75
$ws.Range("N75").Value = -46871.75
$ws.Range("L75").Value = 44999.75
$ws.Range("J78").Value = 44999.75
$ws.Range("H78").Value = 44999.75
$ws.Range("L78").Value = 134999.25
$ws.Range("N78").Value = -144359.25
$ws.Range("H122").Value = 4014.7778
$ws.Range("M122").Value = -5743.2001
$ws.Range("K122").Value = 8193.2001
$ws.Range("I122").Value = 2731.0667
